$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0

$ws.Range("F2").Value = 1309083561.651263
$ws.Range("F3").Value = 1309083561.651263

$ws.Range("D4").Value = 0.098228165751027524
$ws.Range("E4").Value = 0.68607232910464744
$ws.Range("F4").Value = 1309083561.651263

$ws.Range("C5").Value = 0.36
$ws.Range("D5").Value = 0.23769663051268039
$ws.Range("E5").Value = 0.68607232910464744
$ws.Range("F5").Value = 1270172230.762085

$ws.Rows("6:8").Delete()
